$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Statistics")

# Insert a new column before column A, shifting existing data (B:I) right
$ws.Range("A1").EntireColumn.Insert()

# Set the new header and website names
$ws.Range("A1").Value = "Website"
$ws.Range("A2").Value = "blindonion"
$ws.Range("A3").Value = "poop"

# Match the header formatting used by the rest of row 1 (bold, bordered, centered)
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
